$wb = $excel.ActiveWorkbook

# The f0e063cd-b80b-42ab-83e7-261f41b2a5aa.md file has finished its localization
# handback cycle for both locales; generate the handback report:
#  - flip Status from "Ready for handoff" to "Handed back: in sync with en-US"
#  - stamp the Latest Handback DateTime for each locale's row

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Handed back: in sync with en-US"
$zhcn.Range("G3").Value = "2016-03-09 09:54:00"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Handed back: in sync with en-US"
$dede.Range("G3").Value = "2016-03-09 09:54:08"
